# Auto-generated Excel COM-interop script to apply scheduled-runner price updates
# across the Titan_Profits workbook (8 job sheets: ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 3
$ws.Range("H3").Value = 500328.5
$ws.Range("J3").Value = 500328.5
$ws.Range("L3").Value = 500328.5
$ws.Range("N3").Value = -500556.5
# Row 62
$ws.Range("H62").Value = 3593618.5
$ws.Range("I62").Value = 6953619.5
$ws.Range("J62").Value = 9617.200000000001
$ws.Range("K62").Value = 6953619.5
$ws.Range("L62").Value = 9617.200000000001
$ws.Range("M62").Value = -6952995.5
$ws.Range("N62").Value = -10865.2
# Row 65
$ws.Range("H65").Value = 3593618.5
$ws.Range("I65").Value = 6953619.5
$ws.Range("J65").Value = 9617.200000000001
$ws.Range("K65").Value = 34768097.5
$ws.Range("L65").Value = 48086
$ws.Range("M65").Value = -34764977.5
$ws.Range("N65").Value = -54326
# Row 81
$ws.Range("H81").Value = 25000
$ws.Range("J81").Value = 25000
$ws.Range("L81").Value = 25000
$ws.Range("N81").Value = -26996
# Row 84
$ws.Range("H84").Value = 25000
$ws.Range("J84").Value = 25000
$ws.Range("L84").Value = 75000
$ws.Range("N84").Value = -84984
# Row 102
$ws.Range("H102").Value = 500328.5
$ws.Range("J102").Value = 500328.5
$ws.Range("L102").Value = 500328.5
$ws.Range("N102").Value = -506818.5
# Row 129
$ws.Range("H129").Value = 1782.3
$ws.Range("J129").Value = 2378.1428
$ws.Range("L129").Value = 7134.428400000001
$ws.Range("N129").Value = -17134.4284
# Row 132
$ws.Range("H132").Value = 223652.58
$ws.Range("I132").Value = 239135.02
$ws.Range("J132").Value = 26251.5
$ws.Range("K132").Value = 717405.0599999999
$ws.Range("L132").Value = 78754.5
$ws.Range("M132").Value = -714875.0599999999
$ws.Range("N132").Value = -83814.5
# Row 137
$ws.Range("H137").Value = 22728198
$ws.Range("I137").Value = 26316426
$ws.Range("J137").Value = 2753
$ws.Range("K137").Value = 78949278
$ws.Range("L137").Value = 8259
$ws.Range("M137").Value = -78946728
$ws.Range("N137").Value = -13359
# Row 138
$ws.Range("H138").Value = 1342.45
$ws.Range("I138").Value = 723.3889
$ws.Range("J138").Value = 1690.6719
$ws.Range("K138").Value = 2170.1667
$ws.Range("L138").Value = 5072.0157
$ws.Range("M138").Value = 2969.8333
$ws.Range("N138").Value = -15352.0157

$ws = $wb.Worksheets.Item("ARM")
# Row 26
$ws.Range("H26").Value = 3000
$ws.Range("I26").Value = 3000
$ws.Range("K26").Value = 3000
$ws.Range("M26").Value = -2670
# Row 32
$ws.Range("H32").Value = 17513.86
$ws.Range("I32").Value = 3124.0476
$ws.Range("J32").Value = 130833.625
$ws.Range("K32").Value = 3124.0476
$ws.Range("L32").Value = 130833.625
$ws.Range("M32").Value = -2837.0476
$ws.Range("N32").Value = -131407.625
# Row 61
$ws.Range("H61").Value = 2265.9246
$ws.Range("I61").Value = 1526.3256
$ws.Range("J61").Value = 5446.2
$ws.Range("K61").Value = 1526.3256
$ws.Range("L61").Value = 5446.2
$ws.Range("M61").Value = -1314.3256
$ws.Range("N61").Value = -5870.2
# Row 74
$ws.Range("H74").Value = 5309.364
$ws.Range("I74").Value = 1585.75
$ws.Range("J74").Value = 11038
$ws.Range("K74").Value = 1585.75
$ws.Range("L74").Value = 11038
$ws.Range("M74").Value = -711.75
$ws.Range("N74").Value = -12786
# Row 77
$ws.Range("H77").Value = 5309.364
$ws.Range("I77").Value = 1585.75
$ws.Range("J77").Value = 11038
$ws.Range("K77").Value = 7928.75
$ws.Range("L77").Value = 55190
$ws.Range("M77").Value = -3560.75
$ws.Range("N77").Value = -63926
# Row 97
$ws.Range("H97").Value = 4828.3477
$ws.Range("I97").Value = 6035.6113
$ws.Range("J97").Value = 482.2
$ws.Range("K97").Value = 6035.6113
$ws.Range("L97").Value = 482.2
$ws.Range("M97").Value = -5539.6113
$ws.Range("N97").Value = -1474.2
# Row 136
$ws.Range("H136").Value = 2265.9246
$ws.Range("I136").Value = 1526.3256
$ws.Range("J136").Value = 5446.2
$ws.Range("K136").Value = 4578.976799999999
$ws.Range("L136").Value = 16338.6
$ws.Range("M136").Value = -2028.976799999999
$ws.Range("N136").Value = -21438.6

$ws = $wb.Worksheets.Item("BSM")
# Row 86
$ws.Range("H86").Value = 9129.5
$ws.Range("I86").Value = 2243.5557
$ws.Range("K86").Value = 2243.5557
$ws.Range("M86").Value = -1120.5557
# Row 89
$ws.Range("H89").Value = 9129.5
$ws.Range("I89").Value = 2243.5557
$ws.Range("K89").Value = 11217.7785
$ws.Range("M89").Value = -5601.7785
# Row 94
$ws.Range("H94").Value = 1006.1579
$ws.Range("I94").Value = 867.8889
$ws.Range("J94").Value = 3495
$ws.Range("K94").Value = 867.8889
$ws.Range("L94").Value = 3495
$ws.Range("M94").Value = -416.8889
$ws.Range("N94").Value = -4397
# Row 134
$ws.Range("H134").Value = 18520606
$ws.Range("I134").Value = 22223840
$ws.Range("J134").Value = 4434.6665
$ws.Range("K134").Value = 66671520
$ws.Range("L134").Value = 13303.9995
$ws.Range("M134").Value = -66668985
$ws.Range("N134").Value = -18373.9995

$ws = $wb.Worksheets.Item("CRP")
# Row 31
$ws.Range("H31").Value = 1420.5211
$ws.Range("I31").Value = 818.5294
$ws.Range("K31").Value = 818.5294
$ws.Range("M31").Value = -523.5294
# Row 34
$ws.Range("H34").Value = 1420.5211
$ws.Range("I34").Value = 818.5294
$ws.Range("K34").Value = 818.5294
$ws.Range("M34").Value = -616.5294
# Row 74
$ws.Range("H74").Value = 25575.428
$ws.Range("I74").Value = 9000
$ws.Range("J74").Value = 28338
$ws.Range("K74").Value = 9000
$ws.Range("L74").Value = 28338
$ws.Range("N74").Value = -30086
$ws.Range("M74").Value = -8126
# Row 77
$ws.Range("H77").Value = 25575.428
$ws.Range("I77").Value = 9000
$ws.Range("J77").Value = 28338
$ws.Range("K77").Value = 27000
$ws.Range("L77").Value = 85014
$ws.Range("N77").Value = -93750
$ws.Range("M77").Value = -22632
# Row 110
$ws.Range("H110").Value = 41420.8
$ws.Range("J110").Value = 41420.8
$ws.Range("L110").Value = 41420.8
$ws.Range("N110").Value = -49600.8
# Row 132
$ws.Range("H132").Value = 2313.6458
$ws.Range("I132").Value = 1891.6154
$ws.Range("J132").Value = 4142.4443
$ws.Range("K132").Value = 5674.8462
$ws.Range("L132").Value = 12427.3329
$ws.Range("M132").Value = -3144.8462
$ws.Range("N132").Value = -17487.3329
# Row 134
$ws.Range("H134").Value = 2152.561
$ws.Range("I134").Value = 1454.5714
$ws.Range("K134").Value = 4363.7142
$ws.Range("M134").Value = -1828.7142

$ws = $wb.Worksheets.Item("CUL")
# Row 5
$ws.Range("H5").Value = 1410.025
$ws.Range("I5").Value = 1050.0416
$ws.Range("J5").Value = 1950
$ws.Range("K5").Value = 3150.1248
$ws.Range("L5").Value = 5850
$ws.Range("M5").Value = -3038.1248
$ws.Range("N5").Value = -6074
# Row 105
$ws.Range("H105").Value = 6333.3335
$ws.Range("J105").Value = 6333.3335
$ws.Range("L105").Value = 19000.0005
$ws.Range("N105").Value = -24242.0005
# Row 128
$ws.Range("H128").Value = 406239.75
$ws.Range("I128").Value = 406239.75
$ws.Range("K128").Value = 1218719.25
$ws.Range("M128").Value = -1213739.25
# Row 131
$ws.Range("H131").Value = 5557155
$ws.Range("J131").Value = 6412028.5
$ws.Range("L131").Value = 19236085.5
$ws.Range("N131").Value = -19246165.5
# Row 135
$ws.Range("H135").Value = 1410.025
$ws.Range("I135").Value = 1050.0416
$ws.Range("J135").Value = 1950
$ws.Range("K135").Value = 9450.374400000001
$ws.Range("L135").Value = 17550
$ws.Range("M135").Value = -6915.374400000001
$ws.Range("N135").Value = -22620
# Row 140
$ws.Range("H140").Value = 6786.3945
$ws.Range("I140").Value = 8691.654
$ws.Range("J140").Value = 2658.3333
$ws.Range("K140").Value = 26074.962
$ws.Range("L140").Value = 7974.999899999999
$ws.Range("M140").Value = -20894.962
$ws.Range("N140").Value = -18334.9999

$ws = $wb.Worksheets.Item("GSM")
# Row 22
$ws.Range("H22").Value = 52504
$ws.Range("I22").Value = 52504
$ws.Range("K22").Value = 52504
$ws.Range("M22").Value = -51975

$ws = $wb.Worksheets.Item("LTW")
# Row 43
$ws.Range("H43").Value = 13226.286
$ws.Range("J43").Value = 10397.538
$ws.Range("L43").Value = 10397.538
$ws.Range("N43").Value = -10783.538
# Row 55
$ws.Range("H55").Value = 590.06665
$ws.Range("J55").Value = 559.36365
$ws.Range("L55").Value = 559.36365
$ws.Range("N55").Value = -905.36365
# Row 132
$ws.Range("H132").Value = 2127.508
$ws.Range("I132").Value = 1313.0851
$ws.Range("J132").Value = 4519.875
$ws.Range("K132").Value = 3939.2553
$ws.Range("L132").Value = 13559.625
$ws.Range("M132").Value = -1409.2553
$ws.Range("N132").Value = -18619.625

$ws = $wb.Worksheets.Item("WVR")
# Row 113
$ws.Range("H113").Value = 596.3889
$ws.Range("I113").Value = 432.8889
$ws.Range("J113").Value = 759.8889
$ws.Range("K113").Value = 1298.6667
$ws.Range("L113").Value = 2279.6667
$ws.Range("M113").Value = 871.3333
$ws.Range("N113").Value = -6619.6667
# Row 128
$ws.Range("H128").Value = 80000
$ws.Range("J128").Value = 80000
$ws.Range("L128").Value = 80000
$ws.Range("N128").Value = -89960
# Row 132
$ws.Range("H132").Value = 2113.6555
$ws.Range("I132").Value = 2107.6902
$ws.Range("K132").Value = 6323.0706
$ws.Range("M132").Value = -3793.0706
# Row 136
$ws.Range("H136").Value = 35177.535
$ws.Range("I136").Value = 63707.875
$ws.Range("J136").Value = 2571.4285
$ws.Range("K136").Value = 191123.625
$ws.Range("L136").Value = 7714.2855
$ws.Range("M136").Value = -188573.625
$ws.Range("N136").Value = -12814.2855

